$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.499.61"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.521.69"
$ws.Range("E3").Value = "  -0.76%  "
$ws.Range("E4").Value = "  +0.13%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "314.12"
$c.ClearFormats()
$ws.Range("E5").Value = "  +0.76%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "98.76"
$c.ClearFormats()
$ws.Range("E6").Value = "  -1.28%  "
$ws.Range("E7").Value = "  -0.92%  "
$ws.Range("E8").Value = "  +0.16%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.515"
$c.ClearFormats()
$ws.Range("E9").Value = "  -2.42%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "35.15"
$c.ClearFormats()
$ws.Range("E10").Value = "  -2.79%  "
$ws.Range("E11").Value = "  -0.20%  "
$ws.Range("E12").Value = "  +0.90%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "7.17"
$c.ClearFormats()
$ws.Range("D14").Value = "2.912.45"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.525.53"
$ws.Range("E15").Value = "  -2.42%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "15.22"
$c.ClearFormats()
$ws.Range("E16").Value = "  -4.98%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.807"
$c.ClearFormats()
$ws.Range("E17").Value = "  -3.98%  "
$ws.Range("D18").Value = "42.513.42"
$ws.Range("E18").Value = "  +0.16%  "
$ws.Range("E19").Value = "  -2.49%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D20").Value = "0.0₃0938"
$ws.Range("E20").Value = "  -1.21%  "
$ws.Range("B21").Value = "InternetComputer(DFINITY)"
$ws.Range("C21").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "12.11"
$c.ClearFormats()
$ws.Range("E21").Value = "  -1.52%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "68.98"
$c.ClearFormats()
$ws.Range("E22").Value = "  -0.08%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "240.98"
$c.ClearFormats()
$ws.Range("E23").Value = "  -1.23%  "
$ws.Range("E24").Value = "  -1.82%  "
$ws.Range("E25").Value = "  -3.10%  "
$ws.Range("E26").Value = "  +0.03%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "25.44"
$c.ClearFormats()
$ws.Range("E27").Value = "  -3.14%  "
$ws.Range("E28").Value = "  -4.70%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "9.98"
$c.ClearFormats()
$ws.Range("E29").Value = "  -0.97%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "37.72"
$c.ClearFormats()
$ws.Range("E30").Value = "  -6.15%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "5.89"
$c.ClearFormats()
$ws.Range("E31").Value = "  +3.74%  "
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "156.15"
$c.ClearFormats()
$ws.Range("E32").Value = "  -0.73%  "
$ws.Range("E33").Value = "  -3.03%  "
$ws.Range("E34").Value = "  +0.51%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "0.0781"
$c.ClearFormats()
$ws.Range("E35").Value = "  -2.25%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "3.15"
$c.ClearFormats()
$ws.Range("E36").Value = "  -0.95%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.96"
$c.ClearFormats()
$ws.Range("E37").Value = "  -4.61%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "17.55"
$c.ClearFormats()
$ws.Range("E38").Value = "  -3.70%  "
$ws.Range("E39").Value = "  -2.61%  "
$ws.Range("E40").Value = "  -0.66%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "4.20"
$c.ClearFormats()
$ws.Range("E41").Value = "  +0.63%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "21.89"
$c.ClearFormats()
$ws.Range("E42").Value = "  +0.29%  "
$ws.Range("E44").Value = "  -1.02%  "
$ws.Range("D45").Value = "2.006.61"
$ws.Range("E45").Value = "  +2.27%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "3.20"
$c.ClearFormats()
$ws.Range("E46").Value = "  -2.99%  "
$ws.Range("E47").Value = "  +1.11%  "
$ws.Range("D48").Value = "2.765.82"
$ws.Range("E48").Value = "  -1.02%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "78.90"
$c.ClearFormats()
$ws.Range("E49").Value = "  -2.68%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.187"
$c.ClearFormats()
$ws.Range("E50").Value = "  -2.46%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "71.45"
$c.ClearFormats()
$ws.Range("E51").Value = "  -2.21%  "
